$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('B2').Value = '0.258 (0.238 ± 0.020)'
$ws.Range('C2').Value = '00:05:20 (00:26:33 ± 00:10:20)'
$ws.Range('D2').Value = '00:00:00 (00:00:02 ± 00:00:01)'
$ws.Range('B3').Value = '0.229 (0.191 ± 0.021)'
$ws.Range('C3').Value = '00:05:03 (00:06:33 ± 00:01:07)'
$ws.Range('D3').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B4').Value = '0.209 (0.173 ± 0.018)'
$ws.Range('C4').Value = '00:01:00 (00:01:21 ± 00:00:19)'
$ws.Range('D4').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B5').Value = '0.217 (0.171 ± 0.023)'
$ws.Range('C5').Value = '00:05:08 (00:05:16 ± 00:00:04)'
$ws.Range('D5').Value = '00:00:01 (00:00:02 ± 00:00:01)'
$ws.Range('B6').Value = '0.230 (0.198 ± 0.018)'
$ws.Range('C6').Value = '00:04:57 (00:05:00 ± 00:00:03)'
$ws.Range('D6').Value = '00:00:01 (00:00:03 ± 00:00:02)'
$ws.Range('B9').Value = '0.224 (0.157 ± 0.050)'
$ws.Range('C9').Value = '00:05:00 (00:05:07 ± 00:00:15)'
$ws.Range('D9').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B11').Value = '0.180 (0.106 ± 0.044)'
$ws.Range('C11').Value = '00:05:16 (00:05:54 ± 00:00:26)'
$ws.Range('D11').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B13').Value = '0.018 (0.005 ± 0.005)'
$ws.Range('C13').Value = '00:00:07 (00:00:08 ± 00:00:00)'
$ws.Range('D13').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B14').Value = '0.194 (0.144 ± 0.018)'
$ws.Range('C14').Value = '00:02:31 (00:03:12 ± 00:00:18)'
$ws.Range('D14').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B15').Value = '0.246 (0.195 ± 0.020)'
$ws.Range('C15').Value = '00:02:15 (00:04:44 ± 00:00:39)'
$ws.Range('D15').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B16').Value = '0.224 (0.191 ± 0.017)'
$ws.Range('C16').Value = '00:34:35 (00:35:50 ± 00:00:36)'
$ws.Range('D16').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B17').Value = '0.213 (0.184 ± 0.017)'
$ws.Range('C17').Value = '00:05:18 (00:06:48 ± 00:01:27)'
$ws.Range('D17').Value = '00:00:00 (00:00:00 ± 00:00:00)'
